# "SAP number" -> "Vendor number": the header text in column B changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Vendor number"

# Move the selection/active cell to B1 (also clears the previous
# topLeftCell="D1" scroll position since B1 is already in view).
$ws.Range("B1").Select()

# Try to resize the (headless) window to match the saved view state.
# Harmless no-op if the host doesn't expose/persist window geometry.
$win = $excel.ActiveWindow
$win.Width = 51200
$win.Height = 28260
